# Generate Report for Archive
# - Update localization status text from "Ready for handoff" to "In Translation"
#   on every sheet that shows it (Overview!E2:F2, zh-cn!C2, de-de!C2 — they all
#   shared the same string, so all four cells are updated together so the
#   underlying shared string collapses back down to a single entry).
# - Re-fit the two "handoff"/"status" columns that held that text now that the
#   text is shorter (Overview columns E & F, and the Status column C on the
#   zh-cn / de-de detail sheets).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Narrow the now-shorter-text columns to match the refreshed content.
$overview.Range("E1").ColumnWidth = 13.4101845877511
$overview.Range("F1").ColumnWidth = 13.4101845877511
$zhcn.Range("C1").ColumnWidth = 13.4101845877511
$dede.Range("C1").ColumnWidth = 13.4101845877511
